$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 (I0) and J1 (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the existing header formatting (bold, centered, bordered style used
# by B1:H1) by copying the style from H1 onto the new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data cells I2 and J2
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
